# Apply the "Updated symbol list" data refresh described by the diff.
# All target cells are plain text (t="inlineStr") in the workbook, so numeric-
# looking values (prices / percentages) are written with a leading apostrophe
# to force Excel to keep them as text instead of auto-converting to Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''300.90'
$ws.Range("E2").Value = '''-0.97%'
$ws.Range("D3").Value = '''31.42'
$ws.Range("E3").Value = '''-4.13%'
$ws.Range("D4").Value = '''5.161'
$ws.Range("E4").Value = '''-2.54%'
$ws.Range("D5").Value = '''0.07379'
$ws.Range("E5").Value = '''-1.57%'
$ws.Range("D6").Value = '''2.287'
$ws.Range("E6").Value = '''54.99%'
$ws.Range("D7").Value = '''7.899'
$ws.Range("E7").Value = '''0.77%'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '''3.760'
$ws.Range("E8").Value = '''-1.24%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9277'
$ws.Range("E9").Value = '''0.94%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1714'
$ws.Range("E10").Value = '''1.57%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.07579'
$ws.Range("E11").Value = '''-3.46%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.08180'
$ws.Range("E12").Value = '''2.42%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03040'
$ws.Range("E13").Value = '''0.48%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09923'
$ws.Range("E14").Value = '''0.43%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001488'
$ws.Range("E15").Value = '''-0.55%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006080'
$ws.Range("E16").Value = '''-3.59%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.449'
$ws.Range("E17").Value = '''-0.69%'
$ws.Range("D18").Value = '''2.223'
$ws.Range("E18").Value = '''-0.29%'
$ws.Range("D19").Value = '''0.3260'
$ws.Range("E19").Value = '''-2.03%'
$ws.Range("D20").Value = '''0.1337'
$ws.Range("E20").Value = '''0.64%'
$ws.Range("E21").Value = '''3.75%'
$ws.Range("D22").Value = '''0.04649'
$ws.Range("E22").Value = '''0.73%'
$ws.Range("E23").Value = '''-2.36%'
$ws.Range("D24").Value = '''0.001219'
$ws.Range("E24").Value = '''-0.18%'
$ws.Range("D25").Value = '''0.004476'
$ws.Range("E25").Value = '''0.67%'
$ws.Range("D26").Value = '''0.0001299'
$ws.Range("E26").Value = '''-7.15%'
$ws.Range("E27").Value = '''7.58%'
$ws.Range("D39").Value = '''0.01725'
$ws.Range("E39").Value = '''-1.00%'
$ws.Range("D40").Value = '''0.04520'
$ws.Range("E40").Value = '''-0.33%'
$ws.Range("D41").Value = '''0.007130'
$ws.Range("E41").Value = '''-0.62%'
$ws.Range("D42").Value = '''0.1342'
$ws.Range("E42").Value = '''-0.06%'
$ws.Range("D43").Value = '''0.002278'
$ws.Range("E43").Value = '''3.62%'
$ws.Range("D44").Value = '''0.01054'
$ws.Range("E44").Value = '''-16.73%'
$ws.Range("D45").Value = '''0.00006264'
$ws.Range("E45").Value = '''2.40%'
$ws.Range("D46").Value = '''0.006993'
$ws.Range("E46").Value = '''-46.17%'
$ws.Range("D47").Value = '''1.849'
$ws.Range("E47").Value = '''161.10%'
